$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in H1, matching the style used by the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill H2:H15 with 0 (numeric), matching the plain/unstyled data cells
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
